$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 26 (shifts existing rows 26..97 down to 27..98)
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new record's values
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "Femacal de La Calera"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44544
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 100112052
$ws.Range("G26").Value = "Albahaca"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 140
$ws.Range("K26").Value = 4000
$ws.Range("L26").Value = 4500
$ws.Range("M26").Value = 4214
$ws.Range("N26").Value = "$/docena de matas"
$ws.Range("O26").Value = "Provincia de Quillota"
$ws.Range("P26").Value = 702
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = "Hortaliza"
